# Update the two sequence-diagram slides so the "execute()" call shown on the
# TestSessionManager lifeline becomes "execute(model, history)" and widen /
# reposition its label textbox to match (per the target OOXML diff).
#
# Note: Shape.Left/Top/Width/Height are expressed in points (a single-precision
# float) and get multiplied by 12700 to obtain EMUs when the slide XML is
# written back out, so the literals below are the points values that round-trip
# (through that float32 conversion) to the exact target EMU coordinates.

$p = $ppt.ActivePresentation

# Slide 1: "execute()" -> "execute(model, history)"
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(15)
$shp1.Left = 148.1468505859375
$shp1.Top = 195.5961456298828
$shp1.Width = 145.12205505371094
$shp1.Height = 16.964096069335938
$shp1.TextFrame.TextRange.Text = "execute(model, history)"

# Slide 2: "execute()" -> "execute(model, history)"
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(15)
$shp2.Left = 130.57252502441406
$shp2.Top = 195.5961456298828
$shp2.Width = 162.69638061523438
$shp2.Height = 16.964096069335938
$shp2.TextFrame.TextRange.Text = "execute(model, history)"
